$WNS = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-ParaFromXml {
    param([string]$innerXml)

    $d = $word.ActiveDocument
    $lastPara = $d.Paragraphs.Last
    [void]$lastPara.Range.InsertParagraphAfter()

    $d2 = $word.ActiveDocument
    $newPara = $d2.Paragraphs.Last
    [void]$newPara.Range.InsertXML("<w:p $WNS>$innerXml</w:p>")
}

$boldPPr = "<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>"
$boldRun = "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{0}</w:t></w:r>"

# "Model" section -----------------------------------------------------
New-ParaFromXml ($boldPPr + ($boldRun -f "Model"))
New-ParaFromXml $boldPPr
New-ParaFromXml "<w:r><w:t>Cont&#233;m as entidades do banco que v&#227;o ser mapeadas pelo JPA</w:t></w:r>"
New-ParaFromXml ""

# "Repository" section -------------------------------------------------
New-ParaFromXml ($boldPPr + ($boldRun -f "Repository"))
New-ParaFromXml $boldPPr
New-ParaFromXml "<w:r><w:t>Cont&#233;m as interfaces que possuem todas opera&#231;&#245;es b&#225;sicas de CRUD</w:t></w:r>"
New-ParaFromXml ""

# "Controller" section --------------------------------------------------
New-ParaFromXml ($boldPPr + ($boldRun -f "Controller"))
New-ParaFromXml $boldPPr
New-ParaFromXml "<w:r><w:t>Cont&#233;m as classes que fazem o direcionamento de dados para os EndPoints</w:t></w:r>"
New-ParaFromXml ""

# "Service" section -----------------------------------------------------
New-ParaFromXml ($boldPPr + ($boldRun -f "Service"))
New-ParaFromXml ""
New-ParaFromXml "<w:r><w:t xml:space='preserve'>Cont&#233;m as classes que realizam as regras de neg&#243;cio baseado nas interfaces do repositor, exemplo: Validar dados, Calcular valores, Aplicar desconto em vendas. </w:t></w:r>"
